$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")

# New header cell U1 = "area_balcon", formatted like the other header cells (copy format from T1)
$ws.Range("U1").Value = "area_balcon"
$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)

# New data column U2:U6 = 0, formatted with the "0.0" number format + thin border (matches column style)
$cell2 = $ws.Cells.Item(2, 21)
$cell2.NumberFormat = "0.0"
$cell2.Value = 0
$cell2.Borders.LineStyle = 1

$cell2.Copy()
$ws.Range("U3:U6").PasteSpecial(-4122)
$ws.Range("U3").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("U6").Value = 0

Write-Output "done"
